# Integrates skill scoring for candidate assistance:
# Update per-project scoring data and extend the results table with
# additional candidate groups (ids 2-4) and their associated projects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 'LLMGuard'
$ws.Range("C2").Value = 110
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 5

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 'Multi Model Data Analysis for Annotation of Human Activities'
$ws.Range("C3").Value = 83.2
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 83.2
$ws.Range("F3").Value = 5

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 'CloudPhysician''s Vital Extraction Challenge'
$ws.Range("C4").Value = 77.22
$ws.Range("D4").Value = 0.9
$ws.Range("E4").Value = 69.5
$ws.Range("F4").Value = 5

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 'FaceNet Implementation'
$ws.Range("C5").Value = 77.22
$ws.Range("D5").Value = 0.9
$ws.Range("E5").Value = 69.5
$ws.Range("F5").Value = 5

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 'Literature Society IITJ Website'
$ws.Range("C6").Value = 73.92
$ws.Range("D6").Value = 0.9
$ws.Range("E6").Value = 66.53
$ws.Range("F6").Value = 5

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 'SMART SENSING MIDDLEWARE'
$ws.Range("C7").Value = 73.92
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 73.92
$ws.Range("F7").Value = 2

$ws.Range("A8").Value = 2
$ws.Range("B8").Value = 'RAPID'
$ws.Range("C8").Value = 66.75
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 66.75
$ws.Range("F8").Value = 2

$ws.Range("A9").Value = 3
$ws.Range("B9").Value = 'Multi Model Data Analysis for Annotation of Human Activities'
$ws.Range("C9").Value = 107.01
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 100
$ws.Range("F9").Value = 4

$ws.Range("A10").Value = 3
$ws.Range("B10").Value = 'LLMGuard'
$ws.Range("C10").Value = 107.01
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 100
$ws.Range("F10").Value = 4

$ws.Range("A11").Value = 3
$ws.Range("B11").Value = 'Cloudphysician''s Vital Extraction Challenge'
$ws.Range("C11").Value = 100.3
$ws.Range("D11").Value = 0.9
$ws.Range("E11").Value = 90.27
$ws.Range("F11").Value = 4

$ws.Range("A12").Value = 3
$ws.Range("B12").Value = 'Website for the Literature Society of the college'
$ws.Range("C12").Value = 77.22
$ws.Range("D12").Value = 0.9
$ws.Range("E12").Value = 69.5
$ws.Range("F12").Value = 4

$ws.Range("A13").Value = 4
$ws.Range("B13").Value = 'Alcheringa Pass Portal'
$ws.Range("C13").Value = 69.92
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 69.92
$ws.Range("F13").Value = 2

$ws.Range("A14").Value = 4
$ws.Range("B14").Value = 'Video Conferencing Project'
$ws.Range("C14").Value = 66.75
$ws.Range("D14").Value = 0.85
$ws.Range("E14").Value = 56.74
$ws.Range("F14").Value = 2
